$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '53.257.00'
$ws.Range('E2').Value = '  +3.50%  '

$ws.Range('D3').Value = '3.159.26'
$ws.Range('E3').Value = '  +3.87%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.07%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '396.97'
$ws.Range('E5').Value = '  +3.08%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '105.22'
$ws.Range('E6').Value = '  +2.82%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.545'
$ws.Range('E7').Value = '  +0.27%  '

$ws.Range('E8').Value = '  -0.03%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.613'
$ws.Range('E9').Value = '  +5.15%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '38.44'
$ws.Range('E10').Value = '  +4.53%  '

$ws.Range('E11').Value = '  +1.08%  '

$ws.Range('E12').Value = '  +1.00%  '

$ws.Range('D13').Value = '3.650.25'
$ws.Range('E13').Value = '  +3.46%  '

$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '8.05'
$ws.Range('E14').Value = '  +4.18%  '

$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '19.00'
$ws.Range('E15').Value = '  +1.75%  '

$ws.Range('E16').Value = '  +8.28%  '

$ws.Range('D17').Value = '3.150.63'
$ws.Range('E17').Value = '  +3.81%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '10.71'
$ws.Range('E18').Value = '  +1.84%  '

$ws.Range('D19').Value = '53.164.58'
$ws.Range('E19').Value = '  +3.12%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '3.26'
$ws.Range('E20').Value = '  +4.13%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '12.89'
$ws.Range('E21').Value = '  +3.84%  '

$ws.Range('E22').Value = '  +1.38%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '71.25'
$ws.Range('E23').Value = '  +1.44%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '271.33'
$ws.Range('E24').Value = '  +1.02%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '3.21'
$ws.Range('E25').Value = '  +1.71%  '

$ws.Range('E26').Value = '  -1.13%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '27.73'
$ws.Range('E27').Value = '  +2.93%  '

$ws.Range('E28').Value = '  +3.93%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.172'
$ws.Range('E29').Value = '  +0.36%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.999'
$ws.Range('E30').Value = '  -0.02%  '

$ws.Range('E31').Value = '  +1.49%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '10.98'
$ws.Range('E32').Value = '  +7.08%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '37.22'
$ws.Range('E33').Value = '  +7.30%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.0494'
$ws.Range('E34').Value = '  +11.03%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.08'
$ws.Range('E35').Value = '  +0.53%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '50.51'
$ws.Range('E36').Value = '  +0.16%  '

$ws.Range('E37').Value = '  -0.22%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.50'
$ws.Range('E38').Value = '  +4.93%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.77'
$ws.Range('E39').Value = '  +9.70%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '4.09'
$ws.Range('E40').Value = '  +9.61%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.294'
$ws.Range('E41').Value = '  +1.74%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '17.21'
$ws.Range('E42').Value = '  +1.34%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.89'
$ws.Range('E43').Value = '  +1.07%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '129.98'
$ws.Range('E44').Value = '  +1.58%  '

$ws.Range('E45').Value = '  +0.80%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '22.24'
$ws.Range('E46').Value = '  +1.57%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.44'
$ws.Range('E47').Value = '  -1.66%  '

$ws.Range('E48').Value = '  -0.35%  '

$ws.Range('D49').Value = '2.088.96'
$ws.Range('E49').Value = '  +2.99%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0510'
$ws.Range('E50').Value = '  +27.12%  '

$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.928'
$ws.Range('E51').Value = '  +11.57%  '
